# Update existing cell values (odds changed for rows 2, 3, 5, 11)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 2.3
$ws.Cells.Item(2, 9).Value = 3.7
$ws.Cells.Item(2, 10).Value = 3.2
$ws.Cells.Item(2, 12).Value = 4.5
$ws.Cells.Item(2, 23).Value = 5.5
$ws.Cells.Item(2, 24).Value = 9
$ws.Cells.Item(2, 25).Value = 11
$ws.Cells.Item(2, 27).Value = 26
$ws.Cells.Item(2, 34).Value = 7
$ws.Cells.Item(2, 35).Value = 15
$ws.Cells.Item(2, 41).Value = 15
$ws.Cells.Item(2, 52).Value = 81
$ws.Cells.Item(2, 54).Value = 451
$ws.Cells.Item(3, 7).Value = 3
$ws.Cells.Item(3, 9).Value = 2.75
$ws.Cells.Item(3, 12).Value = 3.75
$ws.Cells.Item(3, 24).Value = 12
$ws.Cells.Item(3, 26).Value = 34
$ws.Cells.Item(3, 27).Value = 41
$ws.Cells.Item(3, 29).Value = 4.5
$ws.Cells.Item(3, 35).Value = 11
$ws.Cells.Item(3, 37).Value = 29
$ws.Cells.Item(3, 43).Value = 81
$ws.Cells.Item(3, 46).Value = 1.91
$ws.Cells.Item(3, 49).Value = 4.33
$ws.Cells.Item(3, 52).Value = 67
$ws.Cells.Item(5, 19).Value = 1.83
$ws.Cells.Item(5, 20).Value = 1.98
$ws.Cells.Item(11, 7).Value = 1.67
$ws.Cells.Item(11, 8).Value = 3.5
$ws.Cells.Item(11, 9).Value = 5.5
$ws.Cells.Item(11, 10).Value = 2.38
$ws.Cells.Item(11, 12).Value = 6
$ws.Cells.Item(11, 21).Value = 2.25
$ws.Cells.Item(11, 22).Value = 1.57
$ws.Cells.Item(11, 26).Value = 12
$ws.Cells.Item(11, 29).Value = 7
$ws.Cells.Item(11, 30).Value = 7
$ws.Cells.Item(11, 31).Value = 21
$ws.Cells.Item(11, 41).Value = 9
$ws.Cells.Item(11, 42).Value = 26
$ws.Cells.Item(11, 43).Value = 34
$ws.Cells.Item(11, 44).Value = 67
$ws.Cells.Item(11, 52).Value = 126

# Add the two new rows (13 and 14) for MEXICO - LIGA MX fixtures
# Row 13
$ws.Cells.Item(13, 1).Value = "OtIA6mZO"
$ws.Cells.Item(13, 2).NumberFormat = "@"
$ws.Cells.Item(13, 2).Value = "05/11/2024"
$ws.Cells.Item(13, 3).Value = "22:00"
$ws.Cells.Item(13, 4).Value = "MEXICO - LIGA MX"
$ws.Cells.Item(13, 5).Value = "Club Leon"
$ws.Cells.Item(13, 6).Value = "Mazatlan FC"
$ws.Cells.Item(13, 7).Value = 1.73
$ws.Cells.Item(13, 8).Value = 4
$ws.Cells.Item(13, 9).Value = 4.2
$ws.Cells.Item(13, 10).Value = 2.3
$ws.Cells.Item(13, 11).Value = 2.38
$ws.Cells.Item(13, 12).Value = 4.5
$ws.Cells.Item(13, 13).Value = 1.03
$ws.Cells.Item(13, 14).Value = 15
$ws.Cells.Item(13, 15).Value = 1.2
$ws.Cells.Item(13, 16).Value = 4.33
$ws.Cells.Item(13, 17).Value = 1.67
$ws.Cells.Item(13, 18).Value = 2.15
$ws.Cells.Item(13, 19).Value = 1.33
$ws.Cells.Item(13, 20).Value = 3.25
$ws.Cells.Item(13, 21).Value = 1.67
$ws.Cells.Item(13, 22).Value = 2.1
$ws.Cells.Item(13, 23).Value = 8.5
$ws.Cells.Item(13, 24).Value = 9
$ws.Cells.Item(13, 25).Value = 8.5
$ws.Cells.Item(13, 26).Value = 15
$ws.Cells.Item(13, 27).Value = 13
$ws.Cells.Item(13, 28).Value = 21
$ws.Cells.Item(13, 29).Value = 15
$ws.Cells.Item(13, 30).Value = 7.5
$ws.Cells.Item(13, 31).Value = 15
$ws.Cells.Item(13, 32).Value = 41
$ws.Cells.Item(13, 33).Value = 151
$ws.Cells.Item(13, 34).Value = 15
$ws.Cells.Item(13, 35).Value = 23
$ws.Cells.Item(13, 36).Value = 13
$ws.Cells.Item(13, 37).Value = 41
$ws.Cells.Item(13, 38).Value = 34
$ws.Cells.Item(13, 39).Value = 34
$ws.Cells.Item(13, 40).Value = 4
$ws.Cells.Item(13, 41).Value = 9
$ws.Cells.Item(13, 42).Value = 17
$ws.Cells.Item(13, 43).Value = 29
$ws.Cells.Item(13, 44).Value = 41
$ws.Cells.Item(13, 45).Value = 101
$ws.Cells.Item(13, 46).Value = 3.25
$ws.Cells.Item(13, 47).Value = 7.5
$ws.Cells.Item(13, 48).Value = 51
$ws.Cells.Item(13, 49).Value = 6.5
$ws.Cells.Item(13, 50).Value = 23
$ws.Cells.Item(13, 51).Value = 26
$ws.Cells.Item(13, 52).Value = 67
$ws.Cells.Item(13, 53).Value = 81
$ws.Cells.Item(13, 54).Value = 151
$ws.Cells.Item(13, 55).Value = 501
$ws.Cells.Item(13, 56).Value = 151

# Row 14
$ws.Cells.Item(14, 1).Value = "tYGTM8J5"
$ws.Cells.Item(14, 2).NumberFormat = "@"
$ws.Cells.Item(14, 2).Value = "05/11/2024"
$ws.Cells.Item(14, 3).Value = "22:00"
$ws.Cells.Item(14, 4).Value = "MEXICO - LIGA MX"
$ws.Cells.Item(14, 5).Value = "Santos Laguna"
$ws.Cells.Item(14, 6).Value = "Guadalajara Chivas"
$ws.Cells.Item(14, 7).Value = 5.5
$ws.Cells.Item(14, 8).Value = 4.1
$ws.Cells.Item(14, 9).Value = 1.55
$ws.Cells.Item(14, 10).Value = 5.5
$ws.Cells.Item(14, 11).Value = 2.38
$ws.Cells.Item(14, 12).Value = 2.1
$ws.Cells.Item(14, 13).Value = 1.04
$ws.Cells.Item(14, 14).Value = 13
$ws.Cells.Item(14, 15).Value = 1.22
$ws.Cells.Item(14, 16).Value = 4
$ws.Cells.Item(14, 17).Value = 1.73
$ws.Cells.Item(14, 18).Value = 2.08
$ws.Cells.Item(14, 19).Value = 1.33
$ws.Cells.Item(14, 20).Value = 3.25
$ws.Cells.Item(14, 21).Value = 1.8
$ws.Cells.Item(14, 22).Value = 1.95
$ws.Cells.Item(14, 23).Value = 15
$ws.Cells.Item(14, 24).Value = 29
$ws.Cells.Item(14, 25).Value = 17
$ws.Cells.Item(14, 26).Value = 51
$ws.Cells.Item(14, 27).Value = 41
$ws.Cells.Item(14, 28).Value = 41
$ws.Cells.Item(14, 29).Value = 13
$ws.Cells.Item(14, 30).Value = 8
$ws.Cells.Item(14, 31).Value = 17
$ws.Cells.Item(14, 32).Value = 51
$ws.Cells.Item(14, 33).Value = 251
$ws.Cells.Item(14, 34).Value = 7.5
$ws.Cells.Item(14, 35).Value = 7.5
$ws.Cells.Item(14, 36).Value = 8.5
$ws.Cells.Item(14, 37).Value = 11
$ws.Cells.Item(14, 38).Value = 12
$ws.Cells.Item(14, 39).Value = 23
$ws.Cells.Item(14, 40).Value = 7
$ws.Cells.Item(14, 41).Value = 29
$ws.Cells.Item(14, 42).Value = 34
$ws.Cells.Item(14, 43).Value = 101
$ws.Cells.Item(14, 44).Value = 101
$ws.Cells.Item(14, 45).Value = 201
$ws.Cells.Item(14, 46).Value = 3.25
$ws.Cells.Item(14, 47).Value = 8
$ws.Cells.Item(14, 48).Value = 51
$ws.Cells.Item(14, 49).Value = 3.6
$ws.Cells.Item(14, 50).Value = 7.5
$ws.Cells.Item(14, 51).Value = 17
$ws.Cells.Item(14, 52).Value = 23
$ws.Cells.Item(14, 53).Value = 41
$ws.Cells.Item(14, 54).Value = 126
$ws.Cells.Item(14, 55).Value = 501
$ws.Cells.Item(14, 56).Value = 151
